$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.04"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.396"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06012"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.384"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8109"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9555"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1426"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07421"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03356"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03060"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09424"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.004"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001586"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04809"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005871"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006071"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005049"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009900"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.695"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1342"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03988"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006521"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005229"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005381"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9652"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01917"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
